$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new "turn the pen on" step as row 17, and remove the now
# redundant blank row that used to sit at (old) row 19. Net effect: row
# count / dimension stays the same (A3:D27), but rows 17-19 are
# re-arranged: new r17 = new content, new r18 = old r17, new r19 = old r18.
$ws.Rows("17:17").Insert()
$ws.Range("B17:C17").Merge()
$ws.Range("A17:D17").Borders.LineStyle = 1
$ws.Rows("20:20").Delete()

# --- Fill in the new cell text, in the same order the strings were first
# introduced so the shared-string table comes out in the expected order.
$ws.Range("A24").Value = "Выключить ручку"
$ws.Range("B24").Value = "Ручка выключена"
$ws.Range("D24").Value = "Пройден"

$ws.Range("A17").Value = "3. Включить ручку."
$ws.Range("B17").Value = "Ручка переведена в рабочее состояние"
$ws.Range("D17").Value = "Пройден"

$ws.Range("B7").Value = "Есть механическая ручка."
$ws.Range("B8").Value = "Есть бумага."
$ws.Range("B9").Value = "Есть письменный стол."
$ws.Range("B10").Value = "Есть ступ."
$ws.Range("B11").Value = "Есть линейка."

# --- Match the final cursor/selection state recorded in the sheet.
$ws.Range("B11:C11").Select()
